$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = "Angular routing with kendo grid"
$ws.Range("C24").Value = ":separating controllers when routing using angular ngRoutes"
$ws.Range("A24").Value = 42942

$ws.Range("B22").Select()
